$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 23 with data for 2024-11-01 (serial 45597)
$ws.Range("B23").Value = 0.2826
$ws.Range("C23").Value = 0.27216
$ws.Range("D23").Value = 0.31054
$ws.Range("E23").Value = 0.13474
$ws.Range("F23").Value = 0.04725
$ws.Range("A23").Value = 45597

# Copy the date formatting/style from the row above (A22) onto A23
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)
